# Auto-generated script to apply scheduled-runner value updates to Seraph_Profits workbook
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) - updates market price / profit columns (H-N)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 2639.1  # H111: 2832.3333 -> 2639.1
$ws.Cells.Item(111, 10).Value = 3256.7144  # J111: 3649.5 -> 3256.7144
$ws.Cells.Item(111, 12).Value = 9770.143199999999  # L111: 10948.5 -> 9770.143199999999
$ws.Cells.Item(111, 14).Value = -15904.1432  # N111: -17082.5 -> -15904.1432
$ws.Cells.Item(132, 8).Value = 2388.238  # H132: 2572.3157 -> 2388.238
$ws.Cells.Item(132, 9).Value = 1200.4445  # I132: 1270.5625 -> 1200.4445
$ws.Cells.Item(132, 11).Value = 3601.3335  # K132: 3811.6875 -> 3601.3335
$ws.Cells.Item(132, 13).Value = -1071.3335  # M132: -1281.6875 -> -1071.3335

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 796.8570999999999  # H26: 925.5714 -> 796.8570999999999
$ws.Cells.Item(26, 9).Value = 796.8570999999999  # I26: 925.5714 -> 796.8570999999999
$ws.Cells.Item(26, 11).Value = 796.8570999999999  # K26: 925.5714 -> 796.8570999999999
$ws.Cells.Item(26, 13).Value = -466.8570999999999  # M26: -595.5714 -> -466.8570999999999
$ws.Cells.Item(45, 8).Value = 1708.05  # H45: 1746.7894 -> 1708.05
$ws.Cells.Item(45, 9).Value = 1385.1875  # I45: 1412.7333 -> 1385.1875
$ws.Cells.Item(45, 11).Value = 1385.1875  # K45: 1412.7333 -> 1385.1875
$ws.Cells.Item(45, 13).Value = -1008.1875  # M45: -1035.7333 -> -1008.1875
$ws.Cells.Item(122, 8).Value = 3240.5  # H122: 3342.8462 -> 3240.5
$ws.Cells.Item(122, 10).Value = 4189.143  # J122: 4569 -> 4189.143
$ws.Cells.Item(122, 12).Value = 12567.429  # L122: 13707 -> 12567.429
$ws.Cells.Item(122, 14).Value = -17467.429  # N122: -18607 -> -17467.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 13249.5  # H20: 2861.5715 -> 13249.5
$ws.Cells.Item(20, 9).Value = 11499  # I20: 2724.7144 -> 11499
$ws.Cells.Item(20, 10).Value = 15000  # J20: 2998.4285 -> 15000
$ws.Cells.Item(20, 11).Value = 11499  # K20: 2724.7144 -> 11499
$ws.Cells.Item(20, 12).Value = 15000  # L20: 2998.4285 -> 15000
$ws.Cells.Item(20, 13).Value = -11252  # M20: -2477.7144 -> -11252
$ws.Cells.Item(20, 14).Value = -15494  # N20: -3492.4285 -> -15494
$ws.Cells.Item(86, 8).Value = 2165.6155  # H86: 2219.8462 -> 2165.6155
$ws.Cells.Item(86, 9).Value = 2051.5557  # I86: 2129.889 -> 2051.5557
$ws.Cells.Item(86, 11).Value = 2051.5557  # K86: 2129.889 -> 2051.5557
$ws.Cells.Item(86, 13).Value = -928.5556999999999  # M86: -1006.889 -> -928.5556999999999
$ws.Cells.Item(89, 8).Value = 2165.6155  # H89: 2219.8462 -> 2165.6155
$ws.Cells.Item(89, 9).Value = 2051.5557  # I89: 2129.889 -> 2051.5557
$ws.Cells.Item(89, 11).Value = 10257.7785  # K89: 10649.445 -> 10257.7785
$ws.Cells.Item(89, 13).Value = -4641.7785  # M89: -5033.445 -> -4641.7785
$ws.Cells.Item(95, 8).Value = 10404.6  # H95: 10178 -> 10404.6
$ws.Cells.Item(95, 9).Value = 0  # I95: 3999 -> 0
$ws.Cells.Item(95, 10).Value = 10404.6  # J95: 11207.833 -> 10404.6
$ws.Cells.Item(95, 11).Value = 0  # K95: 3999 -> 0
$ws.Cells.Item(95, 12).Value = 10404.6  # L95: 11207.833 -> 10404.6
$ws.Cells.Item(95, 13).ClearContents()  # M95: was -1253, now empty
$ws.Cells.Item(95, 14).Value = -15896.6  # N95: -16699.833 -> -15896.6
$ws.Cells.Item(105, 8).Value = 8203  # H105: 6487.091 -> 8203
$ws.Cells.Item(105, 9).Value = 6517.857  # I105: 5135.9 -> 6517.857
$ws.Cells.Item(105, 11).Value = 6517.857  # K105: 5135.9 -> 6517.857
$ws.Cells.Item(105, 13).Value = -4770.857  # M105: -3388.9 -> -4770.857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3508.6316  # H31: 3398.3333 -> 3508.6316
$ws.Cells.Item(31, 10).Value = 5853.2856  # J31: 5913.1665 -> 5853.2856
$ws.Cells.Item(31, 12).Value = 5853.2856  # L31: 5913.1665 -> 5853.2856
$ws.Cells.Item(31, 14).Value = -6443.2856  # N31: -6503.1665 -> -6443.2856
$ws.Cells.Item(34, 8).Value = 3508.6316  # H34: 3398.3333 -> 3508.6316
$ws.Cells.Item(34, 10).Value = 5853.2856  # J34: 5913.1665 -> 5853.2856
$ws.Cells.Item(34, 12).Value = 5853.2856  # L34: 5913.1665 -> 5853.2856
$ws.Cells.Item(34, 14).Value = -6257.2856  # N34: -6317.1665 -> -6257.2856
$ws.Cells.Item(35, 8).Value = 387.33334  # H35: 399.85715 -> 387.33334
$ws.Cells.Item(35, 9).Value = 387.33334  # I35: 399.85715 -> 387.33334
$ws.Cells.Item(35, 11).Value = 387.33334  # K35: 399.85715 -> 387.33334
$ws.Cells.Item(35, 13).Value = -93.33334000000002  # M35: -105.85715 -> -93.33334000000002
$ws.Cells.Item(58, 8).Value = 3575.5386  # H58: 3401.5 -> 3575.5386
$ws.Cells.Item(58, 9).Value = 2053.7778  # I58: 1962.3 -> 2053.7778
$ws.Cells.Item(58, 11).Value = 2053.7778  # K58: 1962.3 -> 2053.7778
$ws.Cells.Item(58, 13).Value = -1850.7778  # M58: -1759.3 -> -1850.7778
$ws.Cells.Item(62, 8).Value = 136316.33  # H62: 103237.25 -> 136316.33
$ws.Cells.Item(62, 9).Value = 5000  # I62: 4500 -> 5000
$ws.Cells.Item(62, 11).Value = 5000  # K62: 4500 -> 5000
$ws.Cells.Item(62, 13).Value = -4376  # M62: -3876 -> -4376
$ws.Cells.Item(65, 8).Value = 136316.33  # H65: 103237.25 -> 136316.33
$ws.Cells.Item(65, 9).Value = 5000  # I65: 4500 -> 5000
$ws.Cells.Item(65, 11).Value = 25000  # K65: 22500 -> 25000
$ws.Cells.Item(65, 13).Value = -21880  # M65: -19380 -> -21880
$ws.Cells.Item(99, 8).Value = 15060.583  # H99: 15595.869 -> 15060.583
$ws.Cells.Item(99, 10).Value = 16164.647  # J99: 17003.125 -> 16164.647
$ws.Cells.Item(99, 12).Value = 16164.647  # L99: 17003.125 -> 16164.647
$ws.Cells.Item(99, 14).Value = -19160.647  # N99: -19999.125 -> -19160.647
$ws.Cells.Item(112, 8).Value = 69701.5  # H112: 69702 -> 69701.5
$ws.Cells.Item(112, 10).Value = 69701.5  # J112: 69702 -> 69701.5
$ws.Cells.Item(112, 12).Value = 69701.5  # L112: 69702 -> 69701.5
$ws.Cells.Item(112, 14).Value = -72655.5  # N112: -72656 -> -72655.5
$ws.Cells.Item(126, 8).Value = 15060.583  # H126: 15595.869 -> 15060.583
$ws.Cells.Item(126, 10).Value = 16164.647  # J126: 17003.125 -> 16164.647
$ws.Cells.Item(126, 12).Value = 48493.94100000001  # L126: 51009.375 -> 48493.94100000001
$ws.Cells.Item(126, 14).Value = -53433.94100000001  # N126: -55949.375 -> -53433.94100000001
$ws.Cells.Item(132, 8).Value = 6827.5884  # H132: 7229.524 -> 6827.5884
$ws.Cells.Item(132, 9).Value = 6485.3076  # I132: 6486.4287 -> 6485.3076
$ws.Cells.Item(132, 10).Value = 7940  # J132: 8715.714 -> 7940
$ws.Cells.Item(132, 11).Value = 19455.9228  # K132: 19459.2861 -> 19455.9228
$ws.Cells.Item(132, 12).Value = 23820  # L132: 26147.142 -> 23820
$ws.Cells.Item(132, 13).Value = -16925.9228  # M132: -16929.2861 -> -16925.9228
$ws.Cells.Item(132, 14).Value = -28880  # N132: -31207.142 -> -28880
$ws.Cells.Item(134, 8).Value = 2637.125  # H134: 2887 -> 2637.125
$ws.Cells.Item(134, 9).Value = 2442.5715  # I134: 2519.6 -> 2442.5715
$ws.Cells.Item(134, 10).Value = 3999  # J134: 3499.3333 -> 3999
$ws.Cells.Item(134, 11).Value = 7327.7145  # K134: 7558.799999999999 -> 7327.7145
$ws.Cells.Item(134, 12).Value = 11997  # L134: 10497.9999 -> 11997
$ws.Cells.Item(134, 13).Value = -4792.7145  # M134: -5023.799999999999 -> -4792.7145
$ws.Cells.Item(134, 14).Value = -17067  # N134: -15567.9999 -> -17067
$ws.Cells.Item(136, 8).Value = 3575.5386  # H136: 3401.5 -> 3575.5386
$ws.Cells.Item(136, 9).Value = 2053.7778  # I136: 1962.3 -> 2053.7778
$ws.Cells.Item(136, 11).Value = 6161.3334  # K136: 5886.9 -> 6161.3334
$ws.Cells.Item(136, 13).Value = -3611.3334  # M136: -3336.9 -> -3611.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 1958.7632  # H11: 1873.1666 -> 1958.7632
$ws.Cells.Item(11, 9).Value = 1804.125  # I11: 1848 -> 1804.125
$ws.Cells.Item(11, 10).Value = 2000  # J11: 1999 -> 2000
$ws.Cells.Item(11, 11).Value = 5412.375  # K11: 5544 -> 5412.375
$ws.Cells.Item(11, 12).Value = 6000  # L11: 5997 -> 6000
$ws.Cells.Item(11, 13).Value = -5272.375  # M11: -5404 -> -5272.375
$ws.Cells.Item(11, 14).Value = -6280  # N11: -6277 -> -6280
$ws.Cells.Item(28, 8).Value = 561  # H28: 2090.25 -> 561
$ws.Cells.Item(28, 10).Value = 0  # J28: 2600 -> 0
$ws.Cells.Item(28, 12).Value = 0  # L28: 7800 -> 0
$ws.Cells.Item(28, 14).ClearContents()  # N28: was -8264, now empty
$ws.Cells.Item(34, 8).Value = 1826.3462  # H34: 1958.8636 -> 1826.3462
$ws.Cells.Item(34, 9).Value = 999.375  # I34: 1053.3334 -> 999.375
$ws.Cells.Item(34, 10).Value = 3149.5  # J34: 3899.2856 -> 3149.5
$ws.Cells.Item(34, 11).Value = 2998.125  # K34: 3160.0002 -> 2998.125
$ws.Cells.Item(34, 12).Value = 9448.5  # L34: 11697.8568 -> 9448.5
$ws.Cells.Item(34, 13).Value = -2914.125  # M34: -3076.0002 -> -2914.125
$ws.Cells.Item(34, 14).Value = -9616.5  # N34: -11865.8568 -> -9616.5
$ws.Cells.Item(39, 8).Value = 25475  # H39: 21329.6 -> 25475
$ws.Cells.Item(39, 10).Value = 25475  # J39: 21329.6 -> 25475
$ws.Cells.Item(39, 12).Value = 76425  # L39: 63988.8 -> 76425
$ws.Cells.Item(39, 14).Value = -77013  # N39: -64576.8 -> -77013
$ws.Cells.Item(55, 8).Value = 2599.8572  # H55: 2964.1428 -> 2599.8572
$ws.Cells.Item(55, 10).Value = 2241.5  # J55: 2666.5 -> 2241.5
$ws.Cells.Item(55, 12).Value = 6724.5  # L55: 7999.5 -> 6724.5
$ws.Cells.Item(55, 14).Value = -7078.5  # N55: -8353.5 -> -7078.5
$ws.Cells.Item(106, 8).Value = 0  # H106: 4000 -> 0
$ws.Cells.Item(106, 10).Value = 0  # J106: 4000 -> 0
$ws.Cells.Item(106, 12).Value = 0  # L106: 12000 -> 0
$ws.Cells.Item(106, 14).ClearContents()  # N106: was -13892, now empty
$ws.Cells.Item(118, 8).Value = 709.2  # H118: 690.5 -> 709.2
$ws.Cells.Item(118, 9).Value = 709.2  # I118: 690.5 -> 709.2
$ws.Cells.Item(118, 11).Value = 2127.6  # K118: 2071.5 -> 2127.6
$ws.Cells.Item(118, 13).Value = -884.6000000000004  # M118: -828.5 -> -884.6000000000004
$ws.Cells.Item(120, 8).Value = 11779.6  # H120: 10149.333 -> 11779.6
$ws.Cells.Item(120, 9).Value = 5449.5  # I120: 4299 -> 5449.5
$ws.Cells.Item(120, 11).Value = 16348.5  # K120: 12897 -> 16348.5
$ws.Cells.Item(120, 13).Value = -11510.5  # M120: -8059 -> -11510.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 11831.25  # H70: 8823.947 -> 11831.25
$ws.Cells.Item(70, 9).Value = 10828.333  # I70: 7550.5454 -> 10828.333
$ws.Cells.Item(70, 10).Value = 12834.167  # J70: 10574.875 -> 12834.167
$ws.Cells.Item(70, 11).Value = 10828.333  # K70: 7550.5454 -> 10828.333
$ws.Cells.Item(70, 12).Value = 12834.167  # L70: 10574.875 -> 12834.167
$ws.Cells.Item(70, 13).Value = -10558.333  # M70: -7280.5454 -> -10558.333
$ws.Cells.Item(70, 14).Value = -13374.167  # N70: -11114.875 -> -13374.167
$ws.Cells.Item(73, 8).Value = 11831.25  # H73: 8823.947 -> 11831.25
$ws.Cells.Item(73, 9).Value = 10828.333  # I73: 7550.5454 -> 10828.333
$ws.Cells.Item(73, 10).Value = 12834.167  # J73: 10574.875 -> 12834.167
$ws.Cells.Item(73, 11).Value = 10828.333  # K73: 7550.5454 -> 10828.333
$ws.Cells.Item(73, 12).Value = 12834.167  # L73: 10574.875 -> 12834.167
$ws.Cells.Item(73, 13).Value = -9892.333000000001  # M73: -6614.5454 -> -9892.333000000001
$ws.Cells.Item(73, 14).Value = -14706.167  # N73: -12446.875 -> -14706.167
$ws.Cells.Item(113, 8).Value = 7319.5  # H113: 7323.35 -> 7319.5
$ws.Cells.Item(113, 9).Value = 3289  # I113: 3296.7 -> 3289
$ws.Cells.Item(113, 11).Value = 3289  # K113: 3296.7 -> 3289
$ws.Cells.Item(113, 13).Value = -1119  # M113: -1126.7 -> -1119
$ws.Cells.Item(136, 8).Value = 199999.5  # H136: 199999 -> 199999.5
$ws.Cells.Item(136, 10).Value = 199999.5  # J136: 199999 -> 199999.5
$ws.Cells.Item(136, 12).Value = 599998.5  # L136: 599997 -> 599998.5
$ws.Cells.Item(136, 14).Value = -605098.5  # N136: -605097 -> -605098.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 3443.6  # H4: 3037.8333 -> 3443.6
$ws.Cells.Item(4, 9).Value = 3443.6  # I4: 3037.8333 -> 3443.6
$ws.Cells.Item(4, 11).Value = 3443.6  # K4: 3037.8333 -> 3443.6
$ws.Cells.Item(4, 13).Value = -3330.6  # M4: -2924.8333 -> -3330.6
$ws.Cells.Item(28, 8).Value = 3443.6  # H28: 3037.8333 -> 3443.6
$ws.Cells.Item(28, 9).Value = 3443.6  # I28: 3037.8333 -> 3443.6
$ws.Cells.Item(28, 11).Value = 3443.6  # K28: 3037.8333 -> 3443.6
$ws.Cells.Item(28, 13).Value = -3211.6  # M28: -2805.8333 -> -3211.6
$ws.Cells.Item(32, 8).Value = 0  # H32: 381.5 -> 0
$ws.Cells.Item(32, 9).Value = 0  # I32: 381.5 -> 0
$ws.Cells.Item(32, 11).Value = 0  # K32: 381.5 -> 0
$ws.Cells.Item(32, 13).ClearContents()  # M32: was -64.5, now empty
$ws.Cells.Item(34, 8).Value = 0  # H34: 44000 -> 0
$ws.Cells.Item(34, 9).Value = 0  # I34: 44000 -> 0
$ws.Cells.Item(34, 11).Value = 0  # K34: 44000 -> 0
$ws.Cells.Item(34, 13).ClearContents()  # M34: was -43828, now empty
$ws.Cells.Item(37, 8).Value = 3443.6  # H37: 3037.8333 -> 3443.6
$ws.Cells.Item(37, 9).Value = 3443.6  # I37: 3037.8333 -> 3443.6
$ws.Cells.Item(37, 11).Value = 3443.6  # K37: 3037.8333 -> 3443.6
$ws.Cells.Item(37, 13).Value = -3336.6  # M37: -2930.8333 -> -3336.6
$ws.Cells.Item(40, 8).Value = 3767.5  # H40: 3916.6667 -> 3767.5
$ws.Cells.Item(40, 9).Value = 3767.5  # I40: 3916.6667 -> 3767.5
$ws.Cells.Item(40, 11).Value = 3767.5  # K40: 3916.6667 -> 3767.5
$ws.Cells.Item(40, 13).Value = -3631.5  # M40: -3780.6667 -> -3631.5
$ws.Cells.Item(55, 8).Value = 587  # H55: 603 -> 587
$ws.Cells.Item(55, 9).Value = 587.7778  # I55: 588.3333 -> 587.7778
$ws.Cells.Item(55, 10).Value = 585.25  # J55: 647 -> 585.25
$ws.Cells.Item(55, 11).Value = 587.7778  # K55: 588.3333 -> 587.7778
$ws.Cells.Item(55, 12).Value = 585.25  # L55: 647 -> 585.25
$ws.Cells.Item(55, 13).Value = -414.7778  # M55: -415.3333 -> -414.7778
$ws.Cells.Item(55, 14).Value = -931.25  # N55: -993 -> -931.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 14498.5  # H11: 14499 -> 14498.5
$ws.Cells.Item(11, 10).Value = 16664.666  # J11: 16665.334 -> 16664.666
$ws.Cells.Item(11, 12).Value = 16664.666  # L11: 16665.334 -> 16664.666
$ws.Cells.Item(11, 14).Value = -16948.666  # N11: -16949.334 -> -16948.666
$ws.Cells.Item(33, 8).Value = 29218.666  # H33: 30718.666 -> 29218.666
$ws.Cells.Item(33, 9).Value = 20000  # I33: 29000 -> 20000
$ws.Cells.Item(33, 11).Value = 20000  # K33: 29000 -> 20000
$ws.Cells.Item(33, 13).Value = -19750  # M33: -28750 -> -19750
$ws.Cells.Item(36, 8).Value = 29218.666  # H36: 30718.666 -> 29218.666
$ws.Cells.Item(36, 9).Value = 20000  # I36: 29000 -> 20000
$ws.Cells.Item(36, 11).Value = 20000  # K36: 29000 -> 20000
$ws.Cells.Item(36, 13).Value = -19750  # M36: -28750 -> -19750
$ws.Cells.Item(47, 8).Value = 44749.5  # H47: 44499.5 -> 44749.5
$ws.Cells.Item(47, 9).Value = 0  # I47: 44000 -> 0
$ws.Cells.Item(47, 10).Value = 44749.5  # J47: 44666 -> 44749.5
$ws.Cells.Item(47, 11).Value = 0  # K47: 44000 -> 0
$ws.Cells.Item(47, 12).Value = 44749.5  # L47: 44666 -> 44749.5
$ws.Cells.Item(47, 13).ClearContents()  # M47: was -43428, now empty
$ws.Cells.Item(47, 14).Value = -45893.5  # N47: -45810 -> -45893.5
$ws.Cells.Item(100, 8).Value = 2476.8462  # H100: 2989.9 -> 2476.8462
$ws.Cells.Item(100, 9).Value = 2349.6667  # I100: 2877.3333 -> 2349.6667
$ws.Cells.Item(100, 11).Value = 4699.3334  # K100: 5754.6666 -> 4699.3334
$ws.Cells.Item(100, 13).Value = -4158.3334  # M100: -5213.6666 -> -4158.3334
$ws.Cells.Item(140, 8).Value = 35390  # H140: 59999.5 -> 35390
$ws.Cells.Item(140, 9).Value = 35390  # I140: 0 -> 35390
$ws.Cells.Item(140, 10).Value = 0  # J140: 59999.5 -> 0
$ws.Cells.Item(140, 11).Value = 35390  # K140: 0 -> 35390
$ws.Cells.Item(140, 12).Value = 0  # L140: 59999.5 -> 0
$ws.Cells.Item(140, 13).Value = -30210  # M140: was empty, now -30210
$ws.Cells.Item(140, 14).ClearContents()  # N140: was -70359.5, now empty
